$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 should look like row 2 (same cell formatting), so clone formats first.
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing row 2 (SVM -> SVM_global, refreshed metrics)
$ws.Range("A2").Value = "SVM_global"
$ws.Range("B2").Value = 86.22448979591836
$ws.Range("C2").Value = 92.85714285714286
$ws.Range("D2").Value = 92.85714285714286
$ws.Range("E2").Value = 89.41798941798943
$ws.Range("F2").Value = 0.3153210425937699

# Fill new row 3 (AdaBoostClassifier_global)
$ws.Range("A3").Value = "AdaBoostClassifier_global"
$ws.Range("B3").Value = 86.09271523178806
$ws.Range("C3").Value = 90.90909090909091
$ws.Range("D3").Value = 90.90909090909091
$ws.Range("E3").Value = 88.43537414965988
$ws.Range("F3").Value = 0.659249841068023
